$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values for rows 2-5 (new sensor readings)
$ws.Cells.Item(2, 1).Value2 = 45066.50694444445
$ws.Cells.Item(2, 2).Value2 = 4.411
$ws.Cells.Item(2, 3).Value2 = 1.834
$ws.Cells.Item(2, 4).Value2 = 1.687
$ws.Cells.Item(2, 5).Value2 = 4.96
$ws.Cells.Item(2, 6).Value2 = 6.46
$ws.Cells.Item(2, 7).Value2 = 2.421
$ws.Cells.Item(2, 8).Value2 = 5.343
$ws.Cells.Item(2, 9).Value2 = 2.587
$ws.Cells.Item(2, 10).Value2 = 1.088
$ws.Cells.Item(2, 11).Value2 = 1.625
$ws.Cells.Item(2, 12).Value2 = 2.867
$ws.Cells.Item(2, 13).Value2 = 4.386
$ws.Cells.Item(2, 14).Value2 = 0.5669999999999999
$ws.Cells.Item(2, 15).Value2 = 1.262
$ws.Cells.Item(2, 16).Value2 = 5.05
$ws.Cells.Item(2, 17).Value2 = 0.896
$ws.Cells.Item(2, 18).Value2 = 0.601
$ws.Cells.Item(2, 19).Value2 = 0.104
$ws.Cells.Item(2, 20).Value2 = 31.507
$ws.Cells.Item(2, 21).Value2 = 6.768
$ws.Cells.Item(2, 22).Value2 = 2.777
$ws.Cells.Item(2, 23).Value2 = 3.652
$ws.Cells.Item(2, 24).Value2 = 3.949
$ws.Cells.Item(2, 25).Value2 = 0.176
$ws.Cells.Item(2, 26).Value2 = 1.915
$ws.Cells.Item(2, 27).Value2 = 1.546
$ws.Cells.Item(2, 28).Value2 = 1.916
$ws.Cells.Item(2, 29).Value2 = 5.093
$ws.Cells.Item(2, 30).Value2 = 2.911
$ws.Cells.Item(2, 31).Value2 = 3.891
$ws.Cells.Item(2, 32).Value2 = 4.078
$ws.Cells.Item(2, 33).Value2 = 0.877
$ws.Cells.Item(2, 34).Value2 = 2.065

$ws.Cells.Item(3, 1).Value2 = 45066.51388888889
$ws.Cells.Item(3, 2).Value2 = 14.34
$ws.Cells.Item(3, 3).Value2 = 10.168
$ws.Cells.Item(3, 4).Value2 = 0.953
$ws.Cells.Item(3, 5).Value2 = 29.448
$ws.Cells.Item(3, 6).Value2 = 25.375
$ws.Cells.Item(3, 7).Value2 = 10.912
$ws.Cells.Item(3, 8).Value2 = 36.224
$ws.Cells.Item(3, 9).Value2 = 16.327
$ws.Cells.Item(3, 10).Value2 = 7.314
$ws.Cells.Item(3, 11).Value2 = 10.821
$ws.Cells.Item(3, 12).Value2 = 12.151
$ws.Cells.Item(3, 13).Value2 = 13.314
$ws.Cells.Item(3, 14).Value2 = 3.415
$ws.Cells.Item(3, 15).Value2 = 10.177
$ws.Cells.Item(3, 16).Value2 = 16.27
$ws.Cells.Item(3, 17).Value2 = 8.377000000000001
$ws.Cells.Item(3, 18).Value2 = 0.522
$ws.Cells.Item(3, 19).Value2 = 0.409
$ws.Cells.Item(3, 20).Value2 = 158.19
$ws.Cells.Item(3, 21).Value2 = 30.246
$ws.Cells.Item(3, 22).Value2 = 10.276
$ws.Cells.Item(3, 23).Value2 = 20.176
$ws.Cells.Item(3, 24).Value2 = 11.423
$ws.Cells.Item(3, 25).Value2 = 1.346
$ws.Cells.Item(3, 26).Value2 = 18.2
$ws.Cells.Item(3, 27).Value2 = 8.739000000000001
$ws.Cells.Item(3, 28).Value2 = 7.99
$ws.Cells.Item(3, 29).Value2 = 10.348
$ws.Cells.Item(3, 30).Value2 = 12.769
$ws.Cells.Item(3, 31).Value2 = 1.318
$ws.Cells.Item(3, 32).Value2 = 32.674
$ws.Cells.Item(3, 33).Value2 = 5.474
$ws.Cells.Item(3, 34).Value2 = 12.188

$ws.Cells.Item(4, 1).Value2 = 45066.52083333334
$ws.Cells.Item(4, 2).Value2 = 19.47
$ws.Cells.Item(4, 3).Value2 = 14.236
$ws.Cells.Item(4, 4).Value2 = 0.927
$ws.Cells.Item(4, 5).Value2 = 41.281
$ws.Cells.Item(4, 6).Value2 = 34.782
$ws.Cells.Item(4, 7).Value2 = 15.1
$ws.Cells.Item(4, 8).Value2 = 57.097
$ws.Cells.Item(4, 9).Value2 = 22.951
$ws.Cells.Item(4, 10).Value2 = 10.318
$ws.Cells.Item(4, 11).Value2 = 15.278
$ws.Cells.Item(4, 12).Value2 = 16.771
$ws.Cells.Item(4, 13).Value2 = 18.016
$ws.Cells.Item(4, 14).Value2 = 4.789
$ws.Cells.Item(4, 15).Value2 = 14.526
$ws.Cells.Item(4, 16).Value2 = 21.977
$ws.Cells.Item(4, 17).Value2 = 12.018
$ws.Cells.Item(4, 18).Value2 = 0.463
$ws.Cells.Item(4, 19).Value2 = 0.5639999999999999
$ws.Cells.Item(4, 20).Value2 = 221.328
$ws.Cells.Item(4, 21).Value2 = 41.984
$ws.Cells.Item(4, 22).Value2 = 14.04
$ws.Cells.Item(4, 23).Value2 = 28.246
$ws.Cells.Item(4, 24).Value2 = 15.315
$ws.Cells.Item(4, 25).Value2 = 1.922
$ws.Cells.Item(4, 26).Value2 = 27.9
$ws.Cells.Item(4, 27).Value2 = 12.193
$ws.Cells.Item(4, 28).Value2 = 10.914
$ws.Cells.Item(4, 29).Value2 = 13.396
$ws.Cells.Item(4, 30).Value2 = 17.632
$ws.Cells.Item(4, 31).Value2 = 0.793
$ws.Cells.Item(4, 32).Value2 = 51.957
$ws.Cells.Item(4, 33).Value2 = 7.73
$ws.Cells.Item(4, 34).Value2 = 17.108

$ws.Cells.Item(5, 1).Value2 = 45066.52777777778
$ws.Cells.Item(5, 2).Value2 = 10.27
$ws.Cells.Item(5, 3).Value2 = 7.43
$ws.Cells.Item(5, 4).Value2 = 0.54
$ws.Cells.Item(5, 5).Value2 = 21.57
$ws.Cells.Item(5, 6).Value2 = 18.31
$ws.Cells.Item(5, 7).Value2 = 7.93
$ws.Cells.Item(5, 8).Value2 = 33.92
$ws.Cells.Item(5, 9).Value2 = 11.99
$ws.Cells.Item(5, 10).Value2 = 5.42
$ws.Cells.Item(5, 11).Value2 = 7.94
$ws.Cells.Item(5, 12).Value2 = 8.800000000000001
$ws.Cells.Item(5, 13).Value2 = 9.51
$ws.Cells.Item(5, 14).Value2 = 2.51
$ws.Cells.Item(5, 15).Value2 = 7.5
$ws.Cells.Item(5, 16).Value2 = 11.67
$ws.Cells.Item(5, 17).Value2 = 6.21
$ws.Cells.Item(5, 18).Value2 = 0.33
$ws.Cells.Item(5, 19).Value2 = 0.29
$ws.Cells.Item(5, 20).Value2 = 112.68
$ws.Cells.Item(5, 21).Value2 = 22.16
$ws.Cells.Item(5, 22).Value2 = 7.42
$ws.Cells.Item(5, 23).Value2 = 14.87
$ws.Cells.Item(5, 24).Value2 = 8.119999999999999
$ws.Cells.Item(5, 25).Value2 = 1
$ws.Cells.Item(5, 26).Value2 = 15.91
$ws.Cells.Item(5, 27).Value2 = 6.4
$ws.Cells.Item(5, 28).Value2 = 5.78
$ws.Cells.Item(5, 29).Value2 = 7.18
$ws.Cells.Item(5, 30).Value2 = 9.26
$ws.Cells.Item(5, 31).Value2 = 0.57
$ws.Cells.Item(5, 32).Value2 = 30.99
$ws.Cells.Item(5, 33).Value2 = 4.01
$ws.Cells.Item(5, 34).Value2 = 8.93

# Remove row 6 entirely (dataset now has one fewer row)
$ws.Rows.Item(6).Delete() | Out-Null

# Adjust column widths (custom accuracy formatting tweak)
$ws.Columns.Item(3).ColumnWidth = 7.17  # C -> width 8
$ws.Columns.Item(7).ColumnWidth = 7.17  # G -> width 8
$ws.Columns.Item(10).ColumnWidth = 7.17  # J -> width 8
$ws.Columns.Item(11).ColumnWidth = 7.17  # K -> width 8
$ws.Columns.Item(12).ColumnWidth = 7.17  # L -> width 8
$ws.Columns.Item(13).ColumnWidth = 7.17  # M -> width 8
$ws.Columns.Item(15).ColumnWidth = 7.17  # O -> width 8
$ws.Columns.Item(16).ColumnWidth = 7.17  # P -> width 8
$ws.Columns.Item(17).ColumnWidth = 7.17  # Q -> width 8
$ws.Columns.Item(22).ColumnWidth = 7.17  # V -> width 8
$ws.Columns.Item(24).ColumnWidth = 7.17  # X -> width 8
$ws.Columns.Item(26).ColumnWidth = 6.17  # Z -> width 7
$ws.Columns.Item(27).ColumnWidth = 7.17  # AA -> width 8
$ws.Columns.Item(28).ColumnWidth = 7.17  # AB -> width 8
$ws.Columns.Item(29).ColumnWidth = 7.17  # AC -> width 8
$ws.Columns.Item(30).ColumnWidth = 7.17  # AD -> width 8
$ws.Columns.Item(34).ColumnWidth = 7.17  # AH -> width 8

Write-Output "edit complete"